# Update the cryptos price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped values. Numeric-looking Price values are prefixed
# with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.595.16'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '2.116.57'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  +0.75%  '
$ws.Range("D5").Value = '''336.92'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D7").Value = '''0.5247'
$ws.Range("D8").Value = '''0.4555'
$ws.Range("E8").Value = '  +3.13%  '
$ws.Range("E9").Value = '  +1.81%  '
$ws.Range("D10").Value = '''0.09145'
$ws.Range("E10").Value = '  +2.31%  '
$ws.Range("E11").Value = '  +1.80%  '
$ws.Range("D12").Value = '''24.49'
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").Value = '2.122.53'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '''6.850'
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D15").Value = '''8.133'
$ws.Range("E15").Value = '  +5.82%  '
$ws.Range("D16").Value = '''0.00001182'
$ws.Range("E16").Value = '  +5.34%  '
$ws.Range("D17").Value = '''97.06'
$ws.Range("D18").Value = '''1.010'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D19").Value = '''0.06691'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = '''19.42'
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("D23").Value = '30.654.31'
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").Value = '''12.89'
$ws.Range("E24").Value = '  +4.78%  '
$ws.Range("D25").Value = '''2.356'
$ws.Range("E25").Value = '  +2.00%  '
$ws.Range("D26").Value = '2.360.53'
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").Value = '''164.40'
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("D29").Value = '''2.546'
$ws.Range("E29").Value = '  -0.52%  '
$ws.Range("D30").Value = '''134.68'
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("D31").Value = '''1.211'
$ws.Range("E31").Value = '  +2.03%  '
$ws.Range("D32").Value = '''0.1073'
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").Value = '''1.646'
$ws.Range("E33").Value = '  -0.74%  '
$ws.Range("D34").Value = '''6.368'
$ws.Range("E34").Value = '  +3.53%  '
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("D36").Value = '''10.65'
$ws.Range("E36").Value = '  +5.82%  '
$ws.Range("D37").Value = '''5.873'
$ws.Range("E37").Value = '  +7.28%  '
$ws.Range("D38").Value = '''0.02630'
$ws.Range("E38").Value = '  +2.91%  '
$ws.Range("D39").Value = '''0.06838'
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").Value = '''0.2327'
$ws.Range("E40").Value = '  +3.20%  '
$ws.Range("D41").Value = '''12.59'
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").Value = '''1.255'
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").Value = '''14.80'
$ws.Range("E44").Value = '  +5.93%  '
$ws.Range("D45").Value = '''0.6490'
$ws.Range("E45").Value = '  +2.70%  '
$ws.Range("D46").Value = '''2.313'
$ws.Range("E46").Value = '  +5.36%  '
$ws.Range("D47").Value = '''0.00000000365'
$ws.Range("E47").Value = '  +22.23%  '
$ws.Range("D48").Value = '''3.689'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("D49").Value = '''1.257'
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("E50").Value = '  +2.19%  '
$ws.Range("D51").Value = '''1.184'
$ws.Range("E51").Value = '  -4.34%  '
